# Add "2022-Q4" quarterly data to the workbook:
#  1. Insert a new sheet "2022-Q4" right after "总计" and before "2022-Q3",
#     cloned from "2022-Q3" so the header row / column styles match, then
#     overwrite with the 2022-Q4 figures (and drop the two extra data rows
#     that "2022-Q3" has but "2022-Q4" doesn't).
#  2. Insert a corresponding row into the "总计" summary sheet, shifting the
#     existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the "2022-Q4" sheet by copying "2022-Q3" (keeps header + styles)
# ---------------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcQ3.Copy($srcQ3, $null)

$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# "2022-Q3" has 7 data rows (2..8); "2022-Q4" only has 4 (2..5) - drop the rest.
$newSheet.Rows.Item(6).Resize(3).Delete()

# code, name, fund size, total stock position, position ratio, market value, rank
$q4data = @(
    @("013385", "信澳优势价值混合A", "12.87", "78.75", "3.87", "0.4981", 9),
    @("013393", "信澳价值精选混合A", "3.69", "73.04", "3.12", "0.1151", 10),
    @("013386", "信澳优势价值混合C", "1.30", "78.75", "3.87", "0.0503", 9),
    @("013394", "信澳价值精选混合C", "0.79", "73.04", "3.12", "0.0246", 10)
)

$r = 2
$idx = 0
foreach ($row in $q4data) {
    $newSheet.Cells.Item($r, 1).Value = $idx
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
    $idx = $idx + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a "2022-Q4" row, push the rest down
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$totalData = @(
    @("2022-Q4", 4, 0.69),
    @("2022-Q3", 7, 0.82),
    @("2022-Q2", 8, 0.92),
    @("2022-Q1", 6, 0.33),
    @("2021-Q4", 4, 0.05)
)

# Propagate the row-index column's style (bordered/centered) down to the
# newly-used row 6 before writing values.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$r = 2
$idx = 0
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = $idx
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}
